$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13 : 15684 - 사다리조작
$ws.Range("A13").Value = 15684
$ws.Range("B13").Value = "사다리조작"
$ws.Range("C13").Value = 45905
$ws.Range("C13").NumberFormat = "m/d/yy"

# Row 14 : 14620 - 꽃길
$ws.Range("A14").Value = 14620
$ws.Range("B14").Value = "꽃길"
$ws.Range("C14").Value = 45909
$ws.Range("C14").NumberFormat = "m/d/yy"

# Leave the selection on the last-entered cell, matching the saved view state
$ws.Range("C14").Select()
